# adding heuristica avg error x,y
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coordenadaX (E) and coordenadaY (F) values for rows 2-5
# with heuristic-adjusted (average error corrected) values.
$ws.Range("E2").Value = 806.5705
$ws.Range("F2").Value = 318.68575

$ws.Range("E3").Value = 388.5705
$ws.Range("F3").Value = 489.68575

$ws.Range("E4").Value = 260.5705
$ws.Range("F4").Value = 173.68575

$ws.Range("E5").Value = 630.5705
$ws.Range("F5").Value = 11.31425000000002
